$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "64.618.04"
Set-TextValue $ws.Range("E2") "  +3.59%  "

Set-TextValue $ws.Range("D3") "3.088.87"
Set-TextValue $ws.Range("E3") "  +1.94%  "

Set-TextValue $ws.Range("E4") "  -0.06%  "

Set-TextValue $ws.Range("D5") "558.95"
Set-TextValue $ws.Range("E5") "  +2.97%  "

Set-TextValue $ws.Range("D6") "143.48"
Set-TextValue $ws.Range("E6") "  +7.23%  "

Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  -0.05%  "

Set-TextValue $ws.Range("D8") "3.080.06"
Set-TextValue $ws.Range("E8") "  +1.88%  "

Set-TextValue $ws.Range("E9") "  +0.95%  "

Set-TextValue $ws.Range("D10") "6.35"
Set-TextValue $ws.Range("E10") "  +3.67%  "

Set-TextValue $ws.Range("E11") "  +2.86%  "

Set-TextValue $ws.Range("D12") "0.469"
Set-TextValue $ws.Range("E12") "  +5.01%  "

Set-TextValue $ws.Range("E13") "  +2.58%  "

Set-TextValue $ws.Range("D14") "35.26"
Set-TextValue $ws.Range("E14") "  +2.82%  "

Set-TextValue $ws.Range("D15") "3.600.10"
Set-TextValue $ws.Range("E15") "  +2.25%  "

Set-TextValue $ws.Range("D16") "64.651.06"
Set-TextValue $ws.Range("E16") "  +3.55%  "

Set-TextValue $ws.Range("D17") "3.088.00"
Set-TextValue $ws.Range("E17") "  +1.90%  "

Set-TextValue $ws.Range("E18") "  +0.89%  "

Set-TextValue $ws.Range("D19") "6.77"
Set-TextValue $ws.Range("E19") "  +2.09%  "

Set-TextValue $ws.Range("D20") "477.99"
Set-TextValue $ws.Range("E20") "  -0.25%  "

Set-TextValue $ws.Range("E21") "  +3.33%  "

Set-TextValue $ws.Range("D22") "0.685"
Set-TextValue $ws.Range("E22") "  +1.62%  "

Set-TextValue $ws.Range("D23") "7.56"
Set-TextValue $ws.Range("E23") "  +7.31%  "

Set-TextValue $ws.Range("D24") "13.43"
Set-TextValue $ws.Range("E24") "  +10.65%  "

Set-TextValue $ws.Range("D25") "81.14"
Set-TextValue $ws.Range("E25") "  +0.31%  "

Set-TextValue $ws.Range("E26") "  -0.10%  "

Set-TextValue $ws.Range("D27") "2.78"
Set-TextValue $ws.Range("E27") "  +2.46%  "

Set-TextValue $ws.Range("D28") "8.16"
Set-TextValue $ws.Range("E28") "  +4.78%  "

Set-TextValue $ws.Range("E29") "  +6.25%  "

Set-TextValue $ws.Range("E30") "  +0.19%  "

Set-TextValue $ws.Range("D31") "26.11"
Set-TextValue $ws.Range("E31") "  +1.52%  "

Set-TextValue $ws.Range("E32") "  +2.25%  "

Set-TextValue $ws.Range("D33") "2.47"
Set-TextValue $ws.Range("E33") "  +4.20%  "

Set-TextValue $ws.Range("D34") "5.60"
Set-TextValue $ws.Range("E34") "  -1.21%  "

Set-TextValue $ws.Range("E35") "  +4.41%  "

Set-TextValue $ws.Range("D36") "54.92"
Set-TextValue $ws.Range("E36") "  -0.01%  "

Set-TextValue $ws.Range("D37") "463.07"
Set-TextValue $ws.Range("E37") "  +0.39%  "

Set-TextValue $ws.Range("D38") "0.0832"
Set-TextValue $ws.Range("E38") "  +3.83%  "

Set-TextValue $ws.Range("D39") "0.0408"
Set-TextValue $ws.Range("E39") "  +4.75%  "

Set-TextValue $ws.Range("D40") "2.97"
Set-TextValue $ws.Range("E40") "  +20.05%  "

Set-TextValue $ws.Range("D41") "2.975.46"
Set-TextValue $ws.Range("E41") "  -5.77%  "

Set-TextValue $ws.Range("D42") "8.24"
Set-TextValue $ws.Range("E42") "  +1.67%  "

Set-TextValue $ws.Range("E43") "  -3.08%  "

Set-TextValue $ws.Range("D44") "28.08"
Set-TextValue $ws.Range("E44") "  +5.66%  "

Set-TextValue $ws.Range("E45") "  +5.60%  "

Set-TextValue $ws.Range("E46") "  +0.02%  "

Set-TextValue $ws.Range("D47") "2.15"
Set-TextValue $ws.Range("E47") "  +8.72%  "

Set-TextValue $ws.Range("E48") "  +2.90%  "

Set-TextValue $ws.Range("E49") "  +4.20%  "

Set-TextValue $ws.Range("D50") "116.65"
Set-TextValue $ws.Range("E50") "  +2.09%  "

Set-TextValue $ws.Range("D51") "2.06"
Set-TextValue $ws.Range("E51") "  +1.66%  "
